$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 96.833336
$ws.Range("I6").Value = 79.304344
$ws.Range("K6").Value = 237.913032
$ws.Range("M6").Value = -125.913032
$ws.Range("H101").Value = 3008.4546
$ws.Range("I101").Value = 3398.75
$ws.Range("J101").Value = 2785.4285
$ws.Range("K101").Value = 10196.25
$ws.Range("L101").Value = 8356.2855
$ws.Range("M101").Value = -8574.25
$ws.Range("N101").Value = -11600.2855
$ws.Range("H113").Value = 3083.077
$ws.Range("I113").Value = 3083.077
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3083.077
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170.9229999999998
$ws.Range("N113").ClearContents()
$ws.Range("H125").Value = 1195.5883
$ws.Range("I125").Value = 716.0909
$ws.Range("J125").Value = 2074.6667
$ws.Range("K125").Value = 6444.8181
$ws.Range("L125").Value = 18672.0003
$ws.Range("M125").Value = -3984.8181
$ws.Range("N125").Value = -23592.0003
$ws.Range("H132").Value = 3365.125
$ws.Range("I132").Value = 3312.3333
$ws.Range("J132").Value = 3523.5
$ws.Range("K132").Value = 9936.999899999999
$ws.Range("L132").Value = 10570.5
$ws.Range("M132").Value = -7406.999899999999
$ws.Range("N132").Value = -15630.5
$ws.Range("H137").Value = 2495.9092
$ws.Range("I137").Value = 2328.8333
$ws.Range("K137").Value = 6986.499899999999
$ws.Range("M137").Value = -4436.499899999999
$ws.Range("H138").Value = 3668.392
$ws.Range("J138").Value = 3982.9
$ws.Range("L138").Value = 11948.7
$ws.Range("N138").Value = -22228.7
$ws.Range("H141").Value = 2596.4243
$ws.Range("I141").Value = 2034.7916
$ws.Range("K141").Value = 6104.3748
$ws.Range("M141").Value = -924.3747999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 41670516
$ws.Range("I132").Value = 3936.75
$ws.Range("K132").Value = 11810.25
$ws.Range("M132").Value = -9280.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 14086.6
$ws.Range("I16").Value = 3974.5
$ws.Range("J16").Value = 20828
$ws.Range("K16").Value = 3974.5
$ws.Range("L16").Value = 20828
$ws.Range("M16").Value = -3804.5
$ws.Range("N16").Value = -21168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("H19").Value = 10265873
$ws.Range("I19").Value = 13685190
$ws.Range("J19").Value = 7922.3335
$ws.Range("K19").Value = 13685190
$ws.Range("L19").Value = 7922.3335
$ws.Range("M19").Value = -13685020
$ws.Range("N19").Value = -8262.333500000001
$ws.Range("H24").Value = 10265873
$ws.Range("I24").Value = 13685190
$ws.Range("J24").Value = 7922.3335
$ws.Range("K24").Value = 13685190
$ws.Range("L24").Value = 7922.3335
$ws.Range("M24").Value = -13685020
$ws.Range("N24").Value = -8262.333500000001
$ws.Range("H31").Value = 1841.3077
$ws.Range("I31").Value = 2567.2856
$ws.Range("K31").Value = 2567.2856
$ws.Range("M31").Value = -2272.2856
$ws.Range("H34").Value = 1841.3077
$ws.Range("I34").Value = 2567.2856
$ws.Range("K34").Value = 2567.2856
$ws.Range("M34").Value = -2365.2856
$ws.Range("H58").Value = 3819
$ws.Range("I58").Value = 3669.7
$ws.Range("K58").Value = 3669.7
$ws.Range("M58").Value = -3466.7
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("H134").Value = 1839.4
$ws.Range("I134").Value = 1700
$ws.Range("J134").Value = 1874.25
$ws.Range("K134").Value = 5100
$ws.Range("L134").Value = 5622.75
$ws.Range("M134").Value = -2565
$ws.Range("N134").Value = -10692.75
$ws.Range("H136").Value = 3819
$ws.Range("I136").Value = 3669.7
$ws.Range("K136").Value = 11009.1
$ws.Range("M136").Value = -8459.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 419.14285
$ws.Range("I8").Value = 419.14285
$ws.Range("K8").Value = 1257.42855
$ws.Range("M8").Value = -1118.42855
$ws.Range("H121").Value = 254274.75
$ws.Range("I121").Value = 334033
$ws.Range("J121").Value = 15000
$ws.Range("K121").Value = 1002099
$ws.Range("L121").Value = 45000
$ws.Range("M121").Value = -1000789
$ws.Range("N121").Value = -47620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2850.8572
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 3042.6667
$ws.Range("K80").Value = 1700
$ws.Range("L80").Value = 3042.6667
$ws.Range("M80").Value = -702
$ws.Range("N80").Value = -5038.6667
$ws.Range("H83").Value = 2850.8572
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 3042.6667
$ws.Range("K83").Value = 8500
$ws.Range("L83").Value = 15213.3335
$ws.Range("M83").Value = -3508
$ws.Range("N83").Value = -25197.3335
$ws.Range("H102").Value = 2398.8572
$ws.Range("I102").Value = 2698.75
$ws.Range("K102").Value = 2698.75
$ws.Range("M102").Value = -1076.75
$ws.Range("H113").Value = 3668.2
$ws.Range("I113").Value = 3835.25
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 3835.25
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -1665.25
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 19751.166
$ws.Range("I13").Value = 26250
$ws.Range("J13").Value = 16501.75
$ws.Range("K13").Value = 26250
$ws.Range("L13").Value = 16501.75
$ws.Range("M13").Value = -26110
$ws.Range("N13").Value = -16781.75
$ws.Range("H40").Value = 4336.161
$ws.Range("I40").Value = 3964.96
$ws.Range("J40").Value = 5882.8335
$ws.Range("K40").Value = 3964.96
$ws.Range("L40").Value = 5882.8335
$ws.Range("M40").Value = -3828.96
$ws.Range("N40").Value = -6154.8335
$ws.Range("H61").Value = 2498
$ws.Range("I61").Value = 2498
$ws.Range("K61").Value = 2498
$ws.Range("M61").Value = -2296
$ws.Range("H82").Value = 3097.652
$ws.Range("I82").Value = 2948.1428
$ws.Range("K82").Value = 2948.1428
$ws.Range("M82").Value = -2587.1428
$ws.Range("H85").Value = 3097.652
$ws.Range("I85").Value = 2948.1428
$ws.Range("K85").Value = 2948.1428
$ws.Range("M85").Value = -1700.1428
$ws.Range("H93").Value = 950.4
$ws.Range("I93").Value = 938
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 938
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 310
$ws.Range("N93").Value = -3496
$ws.Range("H113").Value = 2498
$ws.Range("I113").Value = 2498
$ws.Range("K113").Value = 2498
$ws.Range("M113").Value = -328
$ws.Range("H132").Value = 2699.6667
$ws.Range("I132").Value = 2679.6
$ws.Range("K132").Value = 8038.799999999999
$ws.Range("M132").Value = -5508.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19636.5
$ws.Range("J45").Value = 23478.6
$ws.Range("L45").Value = 23478.6
$ws.Range("N45").Value = -24460.6
$ws.Range("H113").Value = 2950
$ws.Range("I113").Value = 2900
$ws.Range("K113").Value = 8700
$ws.Range("M113").Value = -6530
$ws.Range("H122").Value = 2470.5881
$ws.Range("I122").Value = 2491
$ws.Range("J122").Value = 2464.3076
$ws.Range("K122").Value = 7473
$ws.Range("L122").Value = 7392.9228
$ws.Range("M122").Value = -5023
$ws.Range("N122").Value = -12292.9228
$ws.Range("H132").Value = 2711.9285
$ws.Range("I132").Value = 2616.7307
$ws.Range("J132").Value = 3949.5
$ws.Range("K132").Value = 7850.1921
$ws.Range("L132").Value = 11848.5
$ws.Range("M132").Value = -5320.1921
$ws.Range("N132").Value = -16908.5
